$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (renamed/lower-cased, two new columns G/H added)
$ws.Range("A1").Value = "factor"
$ws.Range("B1").Value = "value"
$ws.Range("C1").Value = "dissatisfied_ratio"
$ws.Range("D1").Value = "satisfied_count"
$ws.Range("E1").Value = "dissatisfied_count"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "probability"
$ws.Range("H1").Value = "chi"

# New header cells G1/H1 must use the same header style (bold+border+centered) as the rest of row 1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Sheet shrinks from 32 data+header rows to 31; clear the old trailing row 32
$ws.Range("A32:F32").ClearContents()

# Row 2: 'incident_reopened_flag' / 1
$ws.Range("A2").Value = "incident_reopened_flag"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.7808641975308642
$ws.Range("D2").Value = 213
$ws.Range("E2").Value = 759
$ws.Range("F2").Value = 972
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 4986.45091038875

# Row 3: 'incident_reopened_flag' / 0
$ws.Range("A3").Value = "incident_reopened_flag"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.05341706083836242
$ws.Range("D3").Value = 11607
$ws.Range("E3").Value = 655
$ws.Range("F3").Value = 12262
$ws.Range("G3").Value = [double]"6.544473228086024e-98"
$ws.Range("H3").Value = 4986.45091038875

# Row 4: 'ttr_days_log' / 5
$ws.Range("A4").Value = "ttr_days_log"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 0.4722222222222222
$ws.Range("D4").Value = 152
$ws.Range("E4").Value = 136
$ws.Range("F4").Value = 288
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1077.099187338824

# Row 5: 'ttr_days_log' / 4
$ws.Range("A5").Value = "ttr_days_log"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0.3072060682680152
$ws.Range("D5").Value = 548
$ws.Range("E5").Value = 243
$ws.Range("F5").Value = 791
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1077.099187338824

# Row 6: 'ttr_days_log' / 3
$ws.Range("A6").Value = "ttr_days_log"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 0.1593030491599253
$ws.Range("D6").Value = 2702
$ws.Range("E6").Value = 512
$ws.Range("F6").Value = 3214
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1077.099187338824

# Row 7: 'ttr_days_log' / 2
$ws.Range("A7").Value = "ttr_days_log"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 0.07836644591611479
$ws.Range("D7").Value = 3340
$ws.Range("E7").Value = 284
$ws.Range("F7").Value = 3624
$ws.Range("G7").Value = [double]"4.21433618393792e-09"
$ws.Range("H7").Value = 1077.099187338824

# Row 8: 'ttr_days_log' / 1
$ws.Range("A8").Value = "ttr_days_log"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.04884237234379955
$ws.Range("D8").Value = 2999
$ws.Range("E8").Value = 154
$ws.Range("F8").Value = 3153
$ws.Range("G8").Value = [double]"1.453946301150233e-31"
$ws.Range("H8").Value = 1077.099187338824

# Row 9: 'ttr_days_log' / 0
$ws.Range("A9").Value = "ttr_days_log"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.03927911275415896
$ws.Range("D9").Value = 2079
$ws.Range("E9").Value = 85
$ws.Range("F9").Value = 2164
$ws.Range("G9").Value = [double]"8.346911670100846e-31"
$ws.Range("H9").Value = 1077.099187338824

# Row 10: 'sla_breached' / 1
$ws.Range("A10").Value = "sla_breached"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0.3117106773823192
$ws.Range("D10").Value = 1199
$ws.Range("E10").Value = 543
$ws.Range("F10").Value = 1742
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 879.782667854151

# Row 11: 'sla_breached' / 0
$ws.Range("A11").Value = "sla_breached"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0.0757918552036199
$ws.Range("D11").Value = 10621
$ws.Range("E11").Value = 871
$ws.Range("F11").Value = 11492
$ws.Range("G11").Value = [double]"7.812338348741308e-30"
$ws.Range("H11").Value = 879.782667854151

# Row 12: 'reassignment_count' / 4
$ws.Range("A12").Value = "reassignment_count"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 0.3744075829383886
$ws.Range("D12").Value = 132
$ws.Range("E12").Value = 79
$ws.Range("F12").Value = 211
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 340.4283166328635

# Row 13: 'reassignment_count' / 3
$ws.Range("A13").Value = "reassignment_count"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 0.265017667844523
$ws.Range("D13").Value = 208
$ws.Range("E13").Value = 75
$ws.Range("F13").Value = 283
$ws.Range("G13").Value = 0.9999999999999731
$ws.Range("H13").Value = 340.4283166328635

# Row 14: 'reassignment_count' / 2
$ws.Range("A14").Value = "reassignment_count"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 0.17152466367713
$ws.Range("D14").Value = 739
$ws.Range("E14").Value = 153
$ws.Range("F14").Value = 892
$ws.Range("G14").Value = 0.9999999977984574
$ws.Range("H14").Value = 340.4283166328635

# Row 15: 'reassignment_count' / 1
$ws.Range("A15").Value = "reassignment_count"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 0.1042982277449582
$ws.Range("D15").Value = 4397
$ws.Range("E15").Value = 512
$ws.Range("F15").Value = 4909
$ws.Range("G15").Value = 0.2909728463082871
$ws.Range("H15").Value = 340.4283166328635

# Row 16: 'reassignment_count' / 0
$ws.Range("A16").Value = "reassignment_count"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0.08574722582504683
$ws.Range("D16").Value = 6344
$ws.Range("E16").Value = 595
$ws.Range("F16").Value = 6939
$ws.Range("G16").Value = [double]"2.546909737157177e-09"
$ws.Range("H16").Value = 340.4283166328635

# Row 17: 'close_code' / 'No Resolution Action'
$ws.Range("A17").Value = "close_code"
$ws.Range("B17").Value = "No Resolution Action"
$ws.Range("C17").Value = 0.1696741179639106
$ws.Range("D17").Value = 3083
$ws.Range("E17").Value = 630
$ws.Range("F17").Value = 3713
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 315.57622820523

# Row 18: 'close_code' / 'Information Provided / Training'
$ws.Range("A18").Value = "close_code"
$ws.Range("B18").Value = "Information Provided / Training"
$ws.Range("C18").Value = 0.1089574708317651
$ws.Range("D18").Value = 4735
$ws.Range("E18").Value = 579
$ws.Range("F18").Value = 5314
$ws.Range("G18").Value = 0.7001015665287339
$ws.Range("H18").Value = 315.57622820523

# Row 19: 'close_code' / 'Environmental Restoration'
$ws.Range("A19").Value = "close_code"
$ws.Range("B19").Value = "Environmental Restoration"
$ws.Range("C19").Value = 0.05952380952380952
$ws.Range("D19").Value = 158
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 168
$ws.Range("G19").Value = 0.02475911987982579
$ws.Range("H19").Value = 315.57622820523

# Row 20: 'close_code' / 'Security Modification'
$ws.Range("A20").Value = "close_code"
$ws.Range("B20").Value = "Security Modification"
$ws.Range("C20").Value = 0.05570291777188329
$ws.Range("D20").Value = 712
$ws.Range("E20").Value = 42
$ws.Range("F20").Value = 754
$ws.Range("G20").Value = [double]"5.613066757188756e-07"
$ws.Range("H20").Value = 315.57622820523

# Row 21: 'close_code' / 'Software Correction'
$ws.Range("A21").Value = "close_code"
$ws.Range("B21").Value = "Software Correction"
$ws.Range("C21").Value = 0.05339805825242718
$ws.Range("D21").Value = 390
$ws.Range("E21").Value = 22
$ws.Range("F21").Value = 412
$ws.Range("G21").Value = [double]"9.660615690617236e-05"
$ws.Range("H21").Value = 315.57622820523

# Row 22: 'close_code' / 'Data Correction'
$ws.Range("A22").Value = "close_code"
$ws.Range("B22").Value = "Data Correction"
$ws.Range("C22").Value = 0.0487603305785124
$ws.Range("D22").Value = 2302
$ws.Range("E22").Value = 118
$ws.Range("F22").Value = 2420
$ws.Range("G22").Value = [double]"1.063722824419134e-24"
$ws.Range("H22").Value = 315.57622820523

# Row 23: 'close_code' / 'Reboot / Restart'
$ws.Range("A23").Value = "close_code"
$ws.Range("B23").Value = "Reboot / Restart"
$ws.Range("C23").Value = 0.02869757174392936
$ws.Range("D23").Value = 440
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = 453
$ws.Range("G23").Value = [double]"3.634237019319601e-10"
$ws.Range("H23").Value = 315.57622820523

# Row 24: 'priority_is_4' / 1
$ws.Range("A24").Value = "priority_is_4"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 0.1293998480628007
$ws.Range("D24").Value = 6876
$ws.Range("E24").Value = 1022
$ws.Range("F24").Value = 7898
$ws.Range("G24").Value = 0.9999999998740063
$ws.Range("H24").Value = 103.8263314178949

# Row 25: 'priority_is_4' / 0
$ws.Range("A25").Value = "priority_is_4"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 0.0734632683658171
$ws.Range("D25").Value = 4944
$ws.Range("E25").Value = 392
$ws.Range("F25").Value = 5336
$ws.Range("G25").Value = [double]"5.927722361686574e-17"
$ws.Range("H25").Value = 103.8263314178949

# Row 26: 'caller_is_employee' / 1
$ws.Range("A26").Value = "caller_is_employee"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 0.1245384357166835
$ws.Range("D26").Value = 7824
$ws.Range("E26").Value = 1113
$ws.Range("F26").Value = 8937
$ws.Range("G26").Value = 0.9999999458238565
$ws.Range("H26").Value = 89.71333796111035

# Row 27: 'caller_is_employee' / 0
$ws.Range("A27").Value = "caller_is_employee"
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0.07004887130556202
$ws.Range("D27").Value = 3996
$ws.Range("E27").Value = 301
$ws.Range("F27").Value = 4297
$ws.Range("G27").Value = [double]"8.505400932930321e-17"
$ws.Range("H27").Value = 89.71333796111035

# Row 28: 'incident_has_ka_related_flag' / 0
$ws.Range("A28").Value = "incident_has_ka_related_flag"
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0.1231512666569044
$ws.Range("D28").Value = 5988
$ws.Range("E28").Value = 841
$ws.Range("F28").Value = 6829
$ws.Range("G28").Value = 0.9999912202740207
$ws.Range("H28").Value = 38.95749777616205

# Row 29: 'incident_has_ka_related_flag' / 1
$ws.Range("A29").Value = "incident_has_ka_related_flag"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 0.08946135831381732
$ws.Range("D29").Value = 5832
$ws.Range("E29").Value = 573
$ws.Range("F29").Value = 6405
$ws.Range("G29").Value = [double]"2.185298888607617e-06"
$ws.Range("H29").Value = 38.95749777616205

# Row 30: 'self_service' / 0
$ws.Range("A30").Value = "self_service"
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0.1476846057571965
$ws.Range("D30").Value = 1362
$ws.Range("E30").Value = 236
$ws.Range("F30").Value = 1598
$ws.Range("G30").Value = 0.999999807077853
$ws.Range("H30").Value = 31.27813668428506

# Row 31: 'self_service' / 1
$ws.Range("A31").Value = "self_service"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 0.1012375386730835
$ws.Range("D31").Value = 10458
$ws.Range("E31").Value = 1178
$ws.Range("F31").Value = 11636
$ws.Range("G31").Value = 0.02531359593991387
$ws.Range("H31").Value = 31.27813668428506

